# The workbook "Översikt OKÄNT.xlsx" tracks logged cases in the sheet
# "Avverkningsanmälningar". This update is an automatic refresh of the
# log: the "Förändrad" (last-changed) date in column C moves forward by
# one day (2024-07-22 -> 2024-07-23, serial 45495 -> 45496) for every
# existing data row, and the newest row (row 29, case "A 30076-2024")
# is removed from the log again.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Bump the "Förändrad" date for every existing data row (2-28) from
# 45495 (2024-07-22) to 45496 (2024-07-23).
$ws.Range("C2:C28").Value = 45496

# Remove the last data row (row 29 - case "A 30076-2024"); this also
# shrinks the sheet dimension from A1:Z29 to A1:Z28 and shifts nothing
# else, since it is the final row.
$ws.Rows.Item(29).Delete()

# The row that is now last (28) should no longer carry an explicit
# row height / customHeight flag, matching the target layout.
$ws.Rows.Item(28).AutoFit()
